$d = $word.ActiveDocument

# --- 1. Add the third paragraph after "Este es el segundo cambio" ---
$segundo = $d.Paragraphs(2)
$segundo.Range.InsertParagraphAfter()

$tercero = $d.Paragraphs(3)
$tercero.Range.Text = "Y este el tercero, vamos a ver si podemos recuperar la primera versión"

# --- 2. Move the _GoBack bookmark from the end of paragraph 2 to the end
#        of the new paragraph 3 (it always trails the last edited location). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$tercero = $d.Paragraphs(3)
$endPos = $tercero.Range.End - 1   # just before paragraph 3's own mark

# Inserting the bookmark with a Range collapsed exactly at the end of a
# paragraph's text mis-anchors it, so park a throwaway character there,
# bookmark just in front of it, then remove the throwaway character again.
$d.Range($endPos, $endPos).InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $d.Range($endPos, $endPos))
$d.Range($endPos, $endPos + 1).Delete()
